# Apply "Add data for 2021-11-12" update to the carjacking arrests workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2021-11-04"

# Row 9 (July) - 2021 columns (T/U/V)
$ws.Range("T9").Value = 10
$ws.Range("U9").Value = 139
$ws.Range("V9").Value = 0.0671

# Row 10 (August) - 2021 columns (T/U/V)
$ws.Range("T10").Value = 7
$ws.Range("U10").Value = 153
$ws.Range("V10").Value = 0.0438

# Row 13 (November, partial month) - update label and running totals.
$ws.Range("A13").Value = "November (through 11-04)"
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 10
$ws.Range("I13").Value = 15
$ws.Range("J13").Value = 0.0625
$ws.Range("L13").Value = 14
$ws.Range("M13").Value = 0.125
$ws.Range("O13").Value = 6
$ws.Range("R13").Value = 28
$ws.Range("U13").Value = 24

# Row 14 (Total) - recomputed year totals.
$ws.Range("C14").Value = 229
$ws.Range("D14").Value = 0.1226
$ws.Range("F14").Value = 444
$ws.Range("G14").Value = 0.1048
$ws.Range("I14").Value = 664
$ws.Range("J14").Value = 0.0854
$ws.Range("L14").Value = 563
$ws.Range("M14").Value = 0.1078
$ws.Range("O14").Value = 440
$ws.Range("P14").Value = 0.0984
$ws.Range("R14").Value = 1031
$ws.Range("S14").Value = 0.0498
$ws.Range("T14").Value = 83
$ws.Range("U14").Value = 1385
$ws.Range("V14").Value = 0.0565
